# NEMO.xlsx - add "report_*" flexible-report configuration rows to the
# "configurations" sheet, and restore the various sheet selections /
# active-sheet state that Excel captured when the file was last saved.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. configurations sheet: add new configuration keys/values (rows 3-12)
# ---------------------------------------------------------------------
$configurations = $wb.Worksheets.Item("configurations")

$configurations.Range("A3").Value = "Optimize_DMO_name"

$configurations.Range("A4").Value = "report_title_page"
$configurations.Range("B4").Value = "'True"

$configurations.Range("A5").Value = "report_strategic_challenge"
$configurations.Range("B5").Value = "'True"

$configurations.Range("A6").Value = "report_key_outputs_theme"
$configurations.Range("B6").Value = "'True"

$configurations.Range("A7").Value = "report_decision_makers_options"
$configurations.Range("B7").Value = "'True"

$configurations.Range("A8").Value = "report_scenarios"
$configurations.Range("B8").Value = "'True"

$configurations.Range("A9").Value = "report_fixed_inputs"
$configurations.Range("B9").Value = "'True"

$configurations.Range("A10").Value = "report_dependencies"
$configurations.Range("B10").Value = "'False"

$configurations.Range("A11").Value = "report_weighted_appreciations"
$configurations.Range("B11").Value = "'True"

$configurations.Range("A12").Value = "report_add_optimize"
$configurations.Range("B12").Value = "'False"

# The leading "'" above forces these to be stored as literal text ("True" /
# "False") instead of being auto-coerced to Boolean values. Re-apply the
# default "Normal" style afterwards so the quote-prefix marker doesn't
# change the cells' visible formatting.
$configurations.Range("B4:B12").Style = "Normal"

# Widen the columns to fit the new (longer) configuration keys.
$configurations.Columns.Item(1).ColumnWidth = 25.5
$configurations.Columns.Item(2).ColumnWidth = 12.833333333333334

# ---------------------------------------------------------------------
# 2. Restore per-sheet selections captured in the saved workbook
# ---------------------------------------------------------------------
$decisionMakersOptions = $wb.Worksheets.Item("decision_makers_options")
$decisionMakersOptions.Range("I35").Select() | Out-Null

$dependencies = $wb.Worksheets.Item("dependencies")
$dependencies.Range("C98").Select() | Out-Null

$fixedInputs = $wb.Worksheets.Item("fixed_inputs")
$fixedInputs.Range("A19").Select() | Out-Null

$configurations.Range("B15").Select() | Out-Null

# ---------------------------------------------------------------------
# 3. Make "scenario_weights" the active (selected) tab, as in the saved
#    workbook (was "fixed_inputs" before the edit).
# ---------------------------------------------------------------------
$scenarioWeights = $wb.Worksheets.Item("scenario_weights")
$scenarioWeights.Activate() | Out-Null
